$wb = $excel.ActiveWorkbook

# --- Sheet: Metadata ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "30 Oct 2025, 01:50 PM"

# --- Sheet: Top Gainers ---
$wsGain = $wb.Worksheets.Item("Top Gainers")
$wsGain.Range("B2").Value = "ASALCBR"
$wsGain.Range("C2").Value = 16.4021
$wsGain.Range("D2").Value = 16.6196
$wsGain.Range("E2").Value = 31.4363
$wsGain.Range("B3").Value = "VAIBHAVGBL"
$wsGain.Range("C3").Value = 11.7722
$wsGain.Range("D3").Value = 20.6258
$wsGain.Range("E3").Value = 27.8219
$wsGain.Range("B4").Value = "SALASAR"
$wsGain.Range("C4").Value = 11.269
$wsGain.Range("D4").Value = 16.5957
$wsGain.Range("E4").Value = 23.5626
$wsGain.Range("B5").Value = "KELLTONTEC"
$wsGain.Range("C5").Value = 9.041
$wsGain.Range("D5").Value = 6.1254
$wsGain.Range("E5").Value = -8.5463
$wsGain.Range("B6").Value = "SAGILITY"
$wsGain.Range("C6").Value = 8.7281
$wsGain.Range("D6").Value = 16.3931
$wsGain.Range("E6").Value = 29.6531
$wsGain.Range("B7").Value = "UNIPARTS"
$wsGain.Range("C7").Value = 8.564500000000001
$wsGain.Range("D7").Value = 10.8595
$wsGain.Range("E7").Value = 27.1821
$wsGain.Range("B8").Value = "CHENNPETRO"
$wsGain.Range("C8").Value = 7.5584
$wsGain.Range("D8").Value = 12.7598
$wsGain.Range("E8").Value = 14.9248
$wsGain.Range("B9").Value = "SHREEJISPG"
$wsGain.Range("C9").Value = 7.1643
$wsGain.Range("D9").Value = 11.481
$wsGain.Range("E9").Value = 12.5867
$wsGain.Range("B10").Value = "MARINE"
$wsGain.Range("C10").Value = 6.9461
$wsGain.Range("D10").Value = 3.4555
$wsGain.Range("E10").Value = 16.0889
$wsGain.Range("B11").Value = "EUROPRATIK"
$wsGain.Range("C11").Value = 6.4213
$wsGain.Range("D11").Value = 12.4125
$wsGain.Range("E11").Value = 29.7102
$wsGain.Range("B12").Value = "POLICYBZR"
$wsGain.Range("C12").Value = 6.0752
$wsGain.Range("D12").Value = 8.4475
$wsGain.Range("E12").Value = 7.4089
$wsGain.Range("B13").Value = "BHEL"
$wsGain.Range("C13").Value = 5.7378
$wsGain.Range("D13").Value = 12.3539
$wsGain.Range("E13").Value = 8.815300000000001
$wsGain.Range("B14").Value = "MCLOUD"
$wsGain.Range("C14").Value = 5.7348
$wsGain.Range("D14").Value = 4.7213
$wsGain.Range("E14").Value = -23.3866
$wsGain.Range("B15").Value = "PDSL"
$wsGain.Range("C15").Value = 5.6949
$wsGain.Range("D15").Value = 9.0436
$wsGain.Range("E15").Value = 15.2129
$wsGain.Range("B16").Value = "INDIACEM"
$wsGain.Range("C16").Value = 5.6854
$wsGain.Range("D16").Value = 5.9698
$wsGain.Range("E16").Value = 7.7083
$wsGain.Range("B18").Value = "IVALUE"
$wsGain.Range("C18").Value = 5.2926
$wsGain.Range("D18").Value = 8.6121
$wsGain.Range("E18").Value = 1.2545
$wsGain.Range("B19").Value = "NETWEB"
$wsGain.Range("C19").Value = 5.1093
$wsGain.Range("D19").Value = 10.9655
$wsGain.Range("E19").Value = 13.0519
$wsGain.Range("B24").Value = "VENKEYS"
$wsGain.Range("C24").Value = 4.8503
$wsGain.Range("D24").Value = 5.4763
$wsGain.Range("E24").Value = 3.1821
$wsGain.Range("B25").Value = "VIMTALABS"
$wsGain.Range("C25").Value = 4.8206
$wsGain.Range("D25").Value = 5
$wsGain.Range("E25").Value = -0.1274
$wsGain.Range("B26").Value = "SKYGOLD"
$wsGain.Range("C26").Value = 4.6918
$wsGain.Range("D26").Value = 0.09569999999999999
$wsGain.Range("E26").Value = 39.0124
$wsGain.Range("B27").Value = "MTARTECH"
$wsGain.Range("C27").Value = 4.4964
$wsGain.Range("D27").Value = 8.657
$wsGain.Range("E27").Value = 32.8076
$wsGain.Range("B28").Value = "RAMASTEEL"
$wsGain.Range("C28").Value = 4.4177
$wsGain.Range("D28").Value = 4.3129
$wsGain.Range("E28").Value = 5.9063
$wsGain.Range("B29").Value = "BAJAJHCARE"
$wsGain.Range("C29").Value = 4.4113
$wsGain.Range("D29").Value = 4.96
$wsGain.Range("E29").Value = -1.338
$wsGain.Range("B30").Value = "BLUEDART"
$wsGain.Range("C30").Value = 4.3061
$wsGain.Range("D30").Value = 23.1695
$wsGain.Range("E30").Value = 20.1999
$wsGain.Range("B31").Value = "REDTAPE"
$wsGain.Range("C31").Value = 4.1243
$wsGain.Range("D31").Value = 4.0233
$wsGain.Range("E31").Value = -2.9391
$wsGain.Range("B32").Value = "SHRINGARMS"
$wsGain.Range("C32").Value = 4.0318
$wsGain.Range("D32").Value = 5.2484
$wsGain.Range("E32").Value = 25.5791
$wsGain.Range("B33").Value = "PFOCUS"
$wsGain.Range("C33").Value = 3.9758
$wsGain.Range("D33").Value = 1.2438
$wsGain.Range("E33").Value = 2.7112
$wsGain.Range("B34").Value = "ABREL"
$wsGain.Range("C34").Value = 3.9718
$wsGain.Range("D34").Value = 12.3395
$wsGain.Range("E34").Value = 11.874
$wsGain.Range("B35").Value = "SHANTIGOLD"
$wsGain.Range("C35").Value = 3.9641
$wsGain.Range("D35").Value = 11.2814
$wsGain.Range("E35").Value = 3.8571
$wsGain.Range("B36").Value = "DBL"
$wsGain.Range("C36").Value = 3.8434
$wsGain.Range("D36").Value = 4.9414
$wsGain.Range("E36").Value = 6.0292
$wsGain.Range("B37").Value = "STYLAMIND"
$wsGain.Range("C37").Value = 3.798
$wsGain.Range("D37").Value = 9.594799999999999
$wsGain.Range("E37").Value = 26.7913
$wsGain.Range("B38").Value = "ALICON"
$wsGain.Range("C38").Value = 3.7867
$wsGain.Range("D38").Value = 9.939299999999999
$wsGain.Range("E38").Value = 15.4003
$wsGain.Range("B39").Value = "SAPPHIRE"
$wsGain.Range("C39").Value = 3.7228
$wsGain.Range("D39").Value = 5.5517
$wsGain.Range("E39").Value = 2.8931
$wsGain.Range("B40").Value = "DEEDEV"
$wsGain.Range("C40").Value = 3.6445
$wsGain.Range("D40").Value = -3.2508
$wsGain.Range("E40").Value = -4.0487
$wsGain.Range("B41").Value = "RSYSTEMS"
$wsGain.Range("C41").Value = 3.499
$wsGain.Range("D41").Value = 4.6608
$wsGain.Range("E41").Value = 7.011
$wsGain.Range("B42").Value = "AHLUCONT"
$wsGain.Range("C42").Value = 3.4522
$wsGain.Range("D42").Value = 2.255
$wsGain.Range("E42").Value = -4.999
$wsGain.Range("B43").Value = "CREDITACC"
$wsGain.Range("C43").Value = 3.4053
$wsGain.Range("D43").Value = 2.0387
$wsGain.Range("E43").Value = 7.2643
$wsGain.Range("B44").Value = "CENTRUM"
$wsGain.Range("C44").Value = 3.3333
$wsGain.Range("D44").Value = 1.9432
$wsGain.Range("E44").Value = 1.2771
$wsGain.Range("B46").Value = "TDPOWERSYS"
$wsGain.Range("C46").Value = 3.1906
$wsGain.Range("D46").Value = 7.5545
$wsGain.Range("E46").Value = 16.8655
$wsGain.Range("B47").Value = "BLISSGVS"
$wsGain.Range("C47").Value = 3.1174
$wsGain.Range("D47").Value = 2.4704
$wsGain.Range("E47").Value = 2.8401
$wsGain.Range("B48").Value = "CANBK"
$wsGain.Range("C48").Value = 3.1141
$wsGain.Range("D48").Value = 5.6325
$wsGain.Range("E48").Value = 7.3317
$wsGain.Range("B49").Value = "NEULANDLAB"
$wsGain.Range("C49").Value = 3.0443
$wsGain.Range("D49").Value = -1.324
$wsGain.Range("E49").Value = 8.6957
$wsGain.Range("B50").Value = "OIL"
$wsGain.Range("C50").Value = 3.044
$wsGain.Range("D50").Value = 3.2896
$wsGain.Range("E50").Value = 4.7124
$wsGain.Range("B51").Value = "KMEW"
$wsGain.Range("C51").Value = 3.0394
$wsGain.Range("D51").Value = 5.1613
$wsGain.Range("E51").Value = -0.4291
$wsGain.Range("B52").Value = "VSTIND"
$wsGain.Range("C52").Value = 3.0168
$wsGain.Range("D52").Value = 3.5173
$wsGain.Range("E52").Value = 3.0567
$wsGain.Range("B53").Value = "MRPL"
$wsGain.Range("C53").Value = 2.9545
$wsGain.Range("D53").Value = 12.9516
$wsGain.Range("E53").Value = 23.6012
$wsGain.Range("B54").Value = "RELTD"
$wsGain.Range("C54").Value = 2.92
$wsGain.Range("D54").Value = 9.6881
$wsGain.Range("E54").Value = -1.8027
$wsGain.Range("B55").Value = "ZAGGLE"
$wsGain.Range("C55").Value = 2.9107
$wsGain.Range("D55").Value = 2.8256
$wsGain.Range("E55").Value = 7.7723
$wsGain.Range("B56").Value = "UNIVCABLES"
$wsGain.Range("C56").Value = 2.8983
$wsGain.Range("D56").Value = 3.8481
$wsGain.Range("E56").Value = 3.4661
$wsGain.Range("B57").Value = "GRMOVER"
$wsGain.Range("C57").Value = 2.8979
$wsGain.Range("D57").Value = 3.0721
$wsGain.Range("E57").Value = 18.7404
$wsGain.Range("B58").Value = "REFEX"
$wsGain.Range("C58").Value = 2.8773
$wsGain.Range("D58").Value = 0.3025
$wsGain.Range("E58").Value = 2.2858
$wsGain.Range("B59").Value = "POWERINDIA"
$wsGain.Range("C59").Value = 2.7972
$wsGain.Range("D59").Value = 7.1808
$wsGain.Range("E59").Value = -0.1667
$wsGain.Range("B60").Value = "GMMPFAUDLR"
$wsGain.Range("C60").Value = 2.7909
$wsGain.Range("D60").Value = 7.2456
$wsGain.Range("E60").Value = 19.5621
$wsGain.Range("B61").Value = "JKTYRE"
$wsGain.Range("C61").Value = 2.782
$wsGain.Range("D61").Value = 5.8259
$wsGain.Range("E61").Value = 21.8929
$wsGain.Range("B62").Value = "VOLTAMP"
$wsGain.Range("C62").Value = 2.6843
$wsGain.Range("D62").Value = 2.5634
$wsGain.Range("E62").Value = 2.3012
$wsGain.Range("B63").Value = "MFSL"
$wsGain.Range("C63").Value = 2.6468
$wsGain.Range("D63").Value = 2.7009
$wsGain.Range("E63").Value = -1.0724
$wsGain.Range("B64").Value = "ASHOKA"
$wsGain.Range("C64").Value = 2.5912
$wsGain.Range("D64").Value = 4.1079
$wsGain.Range("E64").Value = 6.7284
$wsGain.Range("B65").Value = "CAMS"
$wsGain.Range("C65").Value = 2.5574
$wsGain.Range("D65").Value = 1.9045
$wsGain.Range("E65").Value = 5.2014
$wsGain.Range("B66").Value = "RGL"
$wsGain.Range("C66").Value = 2.5565
$wsGain.Range("D66").Value = 2.5033
$wsGain.Range("E66").Value = 21.0955
$wsGain.Range("B67").Value = "INOXGREEN"
$wsGain.Range("C67").Value = 2.5518
$wsGain.Range("D67").Value = 10.5473
$wsGain.Range("E67").Value = 33.8869
$wsGain.Range("B68").Value = "WESTLIFE"
$wsGain.Range("C68").Value = 2.5453
$wsGain.Range("D68").Value = 2.5017
$wsGain.Range("E68").Value = -12.4055
$wsGain.Range("B69").Value = "BLS"
$wsGain.Range("C69").Value = 2.5417
$wsGain.Range("D69").Value = -0.487
$wsGain.Range("E69").Value = -1.7375
$wsGain.Range("B70").Value = "CARYSIL"
$wsGain.Range("C70").Value = 2.5078
$wsGain.Range("D70").Value = 1.9929
$wsGain.Range("E70").Value = 10.8671
$wsGain.Range("B71").Value = "FIVESTAR"
$wsGain.Range("C71").Value = 2.4791
$wsGain.Range("D71").Value = 15.4655
$wsGain.Range("E71").Value = 15.5516
$wsGain.Range("B72").Value = "IIFL"
$wsGain.Range("C72").Value = 2.457
$wsGain.Range("D72").Value = 9.4018
$wsGain.Range("E72").Value = 18.5771
$wsGain.Range("B73").Value = "PSPPROJECT"
$wsGain.Range("C73").Value = 2.4446
$wsGain.Range("D73").Value = 16.4598
$wsGain.Range("E73").Value = 22.8627
$wsGain.Range("B74").Value = "JKLAKSHMI"
$wsGain.Range("C74").Value = 2.3908
$wsGain.Range("D74").Value = 4.376
$wsGain.Range("E74").Value = 1.3996
$wsGain.Range("B75").Value = "BPCL"
$wsGain.Range("C75").Value = 2.3844
$wsGain.Range("D75").Value = 7.8529
$wsGain.Range("E75").Value = 4.9315
$wsGain.Range("B76").Value = "WEBELSOLAR"
$wsGain.Range("C76").Value = 2.3546
$wsGain.Range("D76").Value = 2.7507
$wsGain.Range("E76").Value = -1.1091

# --- Sheet: Top Losers ---
$wsLose = $wb.Worksheets.Item("Top Losers")
$wsLose.Range("B2").Value = "IXIGO"
$wsLose.Range("C2").Value = -15.2117
$wsLose.Range("D2").Value = -13.8993
$wsLose.Range("E2").Value = 3.8273
$wsLose.Range("B3").Value = "SHAILY"
$wsLose.Range("C3").Value = -9.6511
$wsLose.Range("D3").Value = -6.218
$wsLose.Range("E3").Value = 6.004
$wsLose.Range("B4").Value = "IDEA"
$wsLose.Range("C4").Value = -6.8376
$wsLose.Range("D4").Value = -9.355499999999999
$wsLose.Range("E4").Value = 7.2571
$wsLose.Range("B5").Value = "KHAICHEM"
$wsLose.Range("C5").Value = -6.2269
$wsLose.Range("D5").Value = -12.529
$wsLose.Range("E5").Value = -10.5831
$wsLose.Range("B10").Value = "NSLNISP"
$wsLose.Range("C10").Value = -4.5041
$wsLose.Range("D10").Value = -3.2592
$wsLose.Range("E10").Value = -4.0571
$wsLose.Range("B11").Value = "LICHSGFIN"
$wsLose.Range("C11").Value = -4.4486
$wsLose.Range("D11").Value = -2.1484
$wsLose.Range("E11").Value = 0.3628
$wsLose.Range("B12").Value = "RAJRATAN"
$wsLose.Range("C12").Value = -4.3
$wsLose.Range("D12").Value = -3.5161
$wsLose.Range("E12").Value = 21.3504
$wsLose.Range("B13").Value = "RAMCOSYS"
$wsLose.Range("C13").Value = -4.2954
$wsLose.Range("D13").Value = 5.421
$wsLose.Range("E13").Value = 23.8307
$wsLose.Range("B14").Value = "KALAMANDIR"
$wsLose.Range("C14").Value = -4.1647
$wsLose.Range("D14").Value = -2.4923
$wsLose.Range("E14").Value = 20.7521
$wsLose.Range("B15").Value = "YATRA"
$wsLose.Range("C15").Value = -4.143
$wsLose.Range("D15").Value = -6.8706
$wsLose.Range("E15").Value = 2.9226
$wsLose.Range("B16").Value = "DRREDDY"
$wsLose.Range("C16").Value = -4.117
$wsLose.Range("D16").Value = -6.5597
$wsLose.Range("E16").Value = -1.9858
$wsLose.Range("B18").Value = "SARDAEN"
$wsLose.Range("C18").Value = -3.8293
$wsLose.Range("D18").Value = -0.3246
$wsLose.Range("E18").Value = -0.3615
$wsLose.Range("B19").Value = "FILATEX"
$wsLose.Range("C19").Value = -3.587
$wsLose.Range("D19").Value = 6.3185
$wsLose.Range("E19").Value = 21.4831
$wsLose.Range("B20").Value = "IDBI"
$wsLose.Range("C20").Value = -3.5129
$wsLose.Range("D20").Value = 4.6843
$wsLose.Range("E20").Value = 7.5585
$wsLose.Range("B21").Value = "TVSHLTD"
$wsLose.Range("C21").Value = -3.4813
$wsLose.Range("D21").Value = -2.2385
$wsLose.Range("E21").Value = 16.0266
$wsLose.Range("B22").Value = "IDEAFORGE"
$wsLose.Range("C22").Value = -3.4745
$wsLose.Range("D22").Value = -2.6754
$wsLose.Range("E22").Value = -4.4125
$wsLose.Range("B23").Value = "DREDGECORP"
$wsLose.Range("C23").Value = -3.4047
$wsLose.Range("D23").Value = 17.6169
$wsLose.Range("E23").Value = 18.3858
$wsLose.Range("B24").Value = "JSL"
$wsLose.Range("C24").Value = -3.347
$wsLose.Range("D24").Value = -3.1482
$wsLose.Range("E24").Value = 5.2931
$wsLose.Range("B25").Value = "VGUARD"
$wsLose.Range("C25").Value = -3.346
$wsLose.Range("D25").Value = -0.6583
$wsLose.Range("E25").Value = -1.5313
$wsLose.Range("B27").Value = "INDUSTOWER"
$wsLose.Range("C27").Value = -3.2804
$wsLose.Range("D27").Value = 1.9361
$wsLose.Range("E27").Value = 7.4803
$wsLose.Range("B28").Value = "SAIL"
$wsLose.Range("C28").Value = -3.2231
$wsLose.Range("D28").Value = 5.0672
$wsLose.Range("E28").Value = 1.1452
$wsLose.Range("B29").Value = "SANDHAR"
$wsLose.Range("C29").Value = -3.1966
$wsLose.Range("D29").Value = 0.4988
$wsLose.Range("E29").Value = 17.77
$wsLose.Range("B30").Value = "UBL"
$wsLose.Range("C30").Value = -3.1571
$wsLose.Range("D30").Value = -2.5641
$wsLose.Range("E30").Value = -1.0329
$wsLose.Range("B32").Value = "NITINSPIN"
$wsLose.Range("C32").Value = -2.986
$wsLose.Range("D32").Value = 0.2314
$wsLose.Range("E32").Value = -0.2762
$wsLose.Range("B33").Value = "GALLANTT"
$wsLose.Range("C33").Value = -2.9458
$wsLose.Range("D33").Value = -1.3082
$wsLose.Range("E33").Value = -20.5692
$wsLose.Range("B34").Value = "GOKULAGRO"
$wsLose.Range("C34").Value = -2.9454
$wsLose.Range("D34").Value = 4.6239
$wsLose.Range("E34").Value = -13.7835
$wsLose.Range("B35").Value = "BIL"
$wsLose.Range("C35").Value = -2.9452
$wsLose.Range("D35").Value = 5.9083
$wsLose.Range("E35").Value = -3.256
$wsLose.Range("B36").Value = "RAYMONDREL"
$wsLose.Range("C36").Value = -2.9308
$wsLose.Range("D36").Value = -4.51
$wsLose.Range("E36").Value = 10.0174
$wsLose.Range("B37").Value = "BCG"
$wsLose.Range("C37").Value = -2.9161
$wsLose.Range("D37").Value = 2.0942
$wsLose.Range("E37").Value = -1.7279
$wsLose.Range("B38").Value = "IEX"
$wsLose.Range("C38").Value = -2.9127
$wsLose.Range("D38").Value = -1.8497
$wsLose.Range("E38").Value = 3.6928
$wsLose.Range("B39").Value = "CGCL"
$wsLose.Range("C39").Value = -2.908
$wsLose.Range("D39").Value = -0.9563
$wsLose.Range("E39").Value = 9.623100000000001
$wsLose.Range("B40").Value = "SURAJEST"
$wsLose.Range("C40").Value = -2.8247
$wsLose.Range("D40").Value = 6.1336
$wsLose.Range("E40").Value = 4.1376
$wsLose.Range("B41").Value = "UTIAMC"
$wsLose.Range("C41").Value = -2.7882
$wsLose.Range("D41").Value = -7.3939
$wsLose.Range("E41").Value = -4.8233
$wsLose.Range("B42").Value = "MANAKCOAT"
$wsLose.Range("C42").Value = -2.7718
$wsLose.Range("D42").Value = -8.616199999999999
$wsLose.Range("E42").Value = 21.3499
$wsLose.Range("B43").Value = "COROMANDEL"
$wsLose.Range("C43").Value = -2.7632
$wsLose.Range("D43").Value = 1.3499
$wsLose.Range("E43").Value = -2.4951
$wsLose.Range("B44").Value = "STYL"
$wsLose.Range("C44").Value = -2.6973
$wsLose.Range("D44").Value = -5.877
$wsLose.Range("E44").Value = -11.4473
$wsLose.Range("B45").Value = "CAMLINFINE"
$wsLose.Range("C45").Value = -2.6426
$wsLose.Range("D45").Value = 0.1379
$wsLose.Range("E45").Value = 0.3995
$wsLose.Range("B46").Value = "HFCL"
$wsLose.Range("C46").Value = -2.6312
$wsLose.Range("D46").Value = -3.3702
$wsLose.Range("E46").Value = 3.0736
$wsLose.Range("B47").Value = "POCL"
$wsLose.Range("C47").Value = -2.601
$wsLose.Range("D47").Value = 2.5566
$wsLose.Range("E47").Value = 23.0554
$wsLose.Range("B48").Value = "BHARATWIRE"
$wsLose.Range("C48").Value = -2.5894
$wsLose.Range("D48").Value = 19.6529
$wsLose.Range("E48").Value = 20.6897
$wsLose.Range("B49").Value = "LXCHEM"
$wsLose.Range("C49").Value = -2.5789
$wsLose.Range("D49").Value = -3.018
$wsLose.Range("E49").Value = -4.1915
$wsLose.Range("B50").Value = "KFINTECH"
$wsLose.Range("C50").Value = -2.5237
$wsLose.Range("D50").Value = -4.2283
$wsLose.Range("E50").Value = 4.6686
$wsLose.Range("B51").Value = "DALMIASUG"
$wsLose.Range("C51").Value = -2.4965
$wsLose.Range("D51").Value = -2.0569
$wsLose.Range("E51").Value = -0.5294
$wsLose.Range("B52").Value = "DCMSRIND"
$wsLose.Range("C52").Value = -2.4913
$wsLose.Range("D52").Value = -1.23
$wsLose.Range("E52").Value = 4.6244
$wsLose.Range("B53").Value = "TVSELECT"
$wsLose.Range("C53").Value = -2.4912
$wsLose.Range("D53").Value = -3.4407
$wsLose.Range("E53").Value = -5.4134
$wsLose.Range("B54").Value = "ARIHANTCAP"
$wsLose.Range("C54").Value = -2.4864
$wsLose.Range("D54").Value = 4.4628
$wsLose.Range("E54").Value = -4.3442
$wsLose.Range("B55").Value = "VINCOFE"
$wsLose.Range("C55").Value = -2.4804
$wsLose.Range("D55").Value = 12.1066
$wsLose.Range("E55").Value = 10.4597
$wsLose.Range("B56").Value = "COSMOFIRST"
$wsLose.Range("C56").Value = -2.4781
$wsLose.Range("D56").Value = -1.4214
$wsLose.Range("E56").Value = -0.2299
$wsLose.Range("B57").Value = "JINDALPHOT"
$wsLose.Range("C57").Value = -2.4376
$wsLose.Range("D57").Value = -2.6533
$wsLose.Range("E57").Value = 19.9855
$wsLose.Range("B58").Value = "VEDL"
$wsLose.Range("C58").Value = -2.4312
$wsLose.Range("D58").Value = 1.6243
$wsLose.Range("E58").Value = 8.1258
$wsLose.Range("B59").Value = "NUVAMA"
$wsLose.Range("C59").Value = -2.4261
$wsLose.Range("D59").Value = 0.0418
$wsLose.Range("E59").Value = 13.862
$wsLose.Range("B60").Value = "VIPIND"
$wsLose.Range("C60").Value = -2.4207
$wsLose.Range("D60").Value = -4.1467
$wsLose.Range("E60").Value = -1.6348
$wsLose.Range("B62").Value = "DCBBANK"
$wsLose.Range("C62").Value = -2.3886
$wsLose.Range("D62").Value = -1.5095
$wsLose.Range("E62").Value = 22.5264
$wsLose.Range("B63").Value = "RPTECH"
$wsLose.Range("C63").Value = -2.3866
$wsLose.Range("D63").Value = -0.6637999999999999
$wsLose.Range("E63").Value = 1.6204
$wsLose.Range("B64").Value = "AEROFLEX"
$wsLose.Range("C64").Value = -2.375
$wsLose.Range("D64").Value = 4.9498
$wsLose.Range("E64").Value = 3.7994
$wsLose.Range("B65").Value = "JMFINANCIL"
$wsLose.Range("C65").Value = -2.372
$wsLose.Range("D65").Value = -2.5599
$wsLose.Range("E65").Value = 4.6337
$wsLose.Range("B66").Value = "IZMO"
$wsLose.Range("C66").Value = -2.3435
$wsLose.Range("D66").Value = -3.4384
$wsLose.Range("E66").Value = -36.3237
$wsLose.Range("B67").Value = "CLEAN"
$wsLose.Range("C67").Value = -2.3344
$wsLose.Range("D67").Value = -3.4749
$wsLose.Range("E67").Value = -6.6994
$wsLose.Range("B68").Value = "PRECWIRE"
$wsLose.Range("C68").Value = -2.3118
$wsLose.Range("D68").Value = 9.7591
$wsLose.Range("E68").Value = 20.436
$wsLose.Range("B69").Value = "UNIONBANK"
$wsLose.Range("C69").Value = -2.3063
$wsLose.Range("D69").Value = 0.6913
$wsLose.Range("E69").Value = 3.0612
$wsLose.Range("B70").Value = "POLYPLEX"
$wsLose.Range("C70").Value = -2.2983
$wsLose.Range("D70").Value = -0.8222
$wsLose.Range("E70").Value = -7.0238
$wsLose.Range("B71").Value = "TTKPRESTIG"
$wsLose.Range("C71").Value = -2.29
$wsLose.Range("D71").Value = 5.528
$wsLose.Range("E71").Value = 7.1395
$wsLose.Range("B72").Value = "EKC"
$wsLose.Range("C72").Value = -2.2866
$wsLose.Range("D72").Value = -3.9547
$wsLose.Range("E72").Value = 0.1374
$wsLose.Range("B73").Value = "TBOTEK"
$wsLose.Range("C73").Value = -2.2779
$wsLose.Range("D73").Value = -5.7697
$wsLose.Range("E73").Value = -1.2655
$wsLose.Range("B74").Value = "TNPL"
$wsLose.Range("C74").Value = -2.2704
$wsLose.Range("D74").Value = 1.4692
$wsLose.Range("E74").Value = -1.8197
$wsLose.Range("B75").Value = "KTKBANK"
$wsLose.Range("C75").Value = -2.2516
$wsLose.Range("D75").Value = -0.0998
$wsLose.Range("E75").Value = 4.6772

# --- Sheet: 1 Month Performance ---
$wsPerf = $wb.Worksheets.Item("1 Month Performance")
$wsPerf.Range("B4").Value = "GMBREW"
$wsPerf.Range("C4").Value = 78.4645
$wsPerf.Range("B5").Value = "IFBAGRO"
$wsPerf.Range("C5").Value = 68.77249999999999
$wsPerf.Range("B7").Value = "MAHASTEEL"
$wsPerf.Range("C7").Value = 54.7315
$wsPerf.Range("B8").Value = "INOXGREEN"
$wsPerf.Range("C8").Value = 53.6434
$wsPerf.Range("B9").Value = "ESSARSHPNG"
$wsPerf.Range("C9").Value = 50.4132
$wsPerf.Range("B10").Value = "NETWEB"
$wsPerf.Range("C10").Value = 45.1184
$wsPerf.Range("B11").Value = "MTARTECH"
$wsPerf.Range("C11").Value = 42.2587
$wsPerf.Range("B15").Value = "SANDUMA"
$wsPerf.Range("C15").Value = 38.0021
$wsPerf.Range("B16").Value = "SHAREINDIA"
$wsPerf.Range("C16").Value = 37.9488
$wsPerf.Range("B17").Value = "BHARATSE"
$wsPerf.Range("C17").Value = 36.2064
$wsPerf.Range("B18").Value = "TVSELECT"
$wsPerf.Range("C18").Value = 35.4622
$wsPerf.Range("B20").Value = "MEGASOFT"
$wsPerf.Range("C20").Value = 35.1224
$wsPerf.Range("B21").Value = "SAMMAANCAP"
$wsPerf.Range("C21").Value = 34.8612
$wsPerf.Range("B24").Value = "INDORAMA"
$wsPerf.Range("C24").Value = 30.8525
$wsPerf.Range("B25").Value = "SOUTHBANK"
$wsPerf.Range("C25").Value = 30.1247
$wsPerf.Range("B26").Value = "ORIENTTECH"
$wsPerf.Range("C26").Value = 30.0154
$wsPerf.Range("B28").Value = "ONMOBILE"
$wsPerf.Range("C28").Value = 28.3888
$wsPerf.Range("B29").Value = "TARACHAND"
$wsPerf.Range("C29").Value = 28.3808
$wsPerf.Range("B30").Value = "MRPL"
$wsPerf.Range("C30").Value = 28.3569
$wsPerf.Range("B31").Value = "ADANIPOWER"
$wsPerf.Range("C31").Value = 27.6693
$wsPerf.Range("B32").Value = "TDPOWERSYS"
$wsPerf.Range("C32").Value = 26.9418
$wsPerf.Range("B33").Value = "SKYGOLD"
$wsPerf.Range("C33").Value = 26.2985
$wsPerf.Range("B34").Value = "HATSUN"
$wsPerf.Range("C34").Value = 25.4153
$wsPerf.Range("B35").Value = "MARINE"
$wsPerf.Range("C35").Value = 25.2394
$wsPerf.Range("B36").Value = "CARTRADE"
$wsPerf.Range("C36").Value = 25.159
$wsPerf.Range("B37").Value = "EMKAY"
$wsPerf.Range("C37").Value = 25.1422
$wsPerf.Range("B39").Value = "UNIPARTS"
$wsPerf.Range("C39").Value = 24.6351
$wsPerf.Range("B40").Value = "AVALON"
$wsPerf.Range("C40").Value = 24.4657
$wsPerf.Range("B41").Value = "ATHERENERG"
$wsPerf.Range("C41").Value = 24.4562
$wsPerf.Range("B43").Value = "SAGILITY"
$wsPerf.Range("C43").Value = 23.653
$wsPerf.Range("B44").Value = "AUBANK"
$wsPerf.Range("C44").Value = 23.4252
$wsPerf.Range("B46").Value = "RAMCOSYS"
$wsPerf.Range("C46").Value = 22.8236
$wsPerf.Range("B47").Value = "INDIANB"
$wsPerf.Range("C47").Value = 22.5876
$wsPerf.Range("B48").Value = "DCBBANK"
$wsPerf.Range("C48").Value = 21.9778
$wsPerf.Range("B49").Value = "GUJTHEM"
$wsPerf.Range("C49").Value = 21.823
$wsPerf.Range("B50").Value = "RBLBANK"
$wsPerf.Range("C50").Value = 21.7733
$wsPerf.Range("B51").Value = "GRMOVER"
$wsPerf.Range("C51").Value = 21.7429
$wsPerf.Range("B55").Value = "INDRAMEDCO"
$wsPerf.Range("C55").Value = 21.3764
$wsPerf.Range("B56").Value = "SCI"
$wsPerf.Range("C56").Value = 21.2821
$wsPerf.Range("B57").Value = "STYLAMIND"
$wsPerf.Range("C57").Value = 21.2051
$wsPerf.Range("B58").Value = "BHAGERIA"
$wsPerf.Range("C58").Value = 21.1237
$wsPerf.Range("B59").Value = "SKMEGGPROD"
$wsPerf.Range("C59").Value = 20.7435
$wsPerf.Range("B62").Value = "BHARATWIRE"
$wsPerf.Range("C62").Value = 20.1761
$wsPerf.Range("B63").Value = "LORDSCHLO"
$wsPerf.Range("C63").Value = 20.1623
$wsPerf.Range("B64").Value = "HINDCOPPER"
$wsPerf.Range("C64").Value = 20.038
$wsPerf.Range("B65").Value = "IIFL"
$wsPerf.Range("C65").Value = 19.9575
$wsPerf.Range("B66").Value = "FEDERALBNK"
$wsPerf.Range("C66").Value = 19.775
$wsPerf.Range("B67").Value = "ETHOSLTD"
$wsPerf.Range("C67").Value = 19.7071
$wsPerf.Range("B68").Value = "ASALCBR"
$wsPerf.Range("C68").Value = 19.6285
$wsPerf.Range("B69").Value = "BLUEDART"
$wsPerf.Range("C69").Value = 19.5813
$wsPerf.Range("B70").Value = "WHEELS"
$wsPerf.Range("C70").Value = 19.3785
$wsPerf.Range("B71").Value = "SHRIRAMFIN"
$wsPerf.Range("C71").Value = 19.3334
$wsPerf.Range("B72").Value = "PRECWIRE"
$wsPerf.Range("C72").Value = 19.2661
$wsPerf.Range("B75").Value = "MCX"
$wsPerf.Range("C75").Value = 18.7496
$wsPerf.Range("B76").Value = "THOMASCOTT"
$wsPerf.Range("C76").Value = 18.589

# --- Sheet: distance from Dma50 ---
$wsDma = $wb.Worksheets.Item("distance from Dma50")
$wsDma.Range("C2").Value = 9.692
$wsDma.Range("C3").Value = 7.227
$wsDma.Range("C4").Value = 5.9386
$wsDma.Range("C5").Value = 5.0833
$wsDma.Range("C6").Value = 5.0402
$wsDma.Range("C7").Value = 4.7979
$wsDma.Range("C8").Value = 4.4345
$wsDma.Range("C9").Value = 4.2893
$wsDma.Range("C10").Value = 3.7814
$wsDma.Range("C11").Value = 3.4482
$wsDma.Range("C12").Value = 3.2891
$wsDma.Range("C13").Value = 3.2645
$wsDma.Range("C14").Value = 2.9477
$wsDma.Range("C15").Value = 2.9144
$wsDma.Range("C16").Value = 2.8309
$wsDma.Range("C17").Value = 2.6862
$wsDma.Range("C18").Value = 2.509
$wsDma.Range("C19").Value = 2.4443
$wsDma.Range("C20").Value = 2.1931
$wsDma.Range("C21").Value = 2.1566
$wsDma.Range("C22").Value = 1.3508
$wsDma.Range("C23").Value = 1.2616
$wsDma.Range("C24").Value = 1.2277
$wsDma.Range("C25").Value = 0.9557
$wsDma.Range("C26").Value = 0.9157999999999999
$wsDma.Range("C27").Value = 0.8655
$wsDma.Range("C28").Value = 0.5409
$wsDma.Range("C29").Value = 0.197
$wsDma.Range("C30").Value = -2.139

Write-Host "Edit applied successfully"
